# Stricter separation of departement specific pre mid and post mid courses
#
# This script:
#  1. Updates the "Allocated Rooms" lists on Basket_Course_Allocations
#     (rows 9-17, column C) so each elective basket course keeps only its
#     own pre-mid rooms (no more overlap with other baskets' rooms).
#  2. Re-points the individual room assignment (column M = "room") on
#     Classroom_Allocation for every affected schedule row, and refreshes
#     the dependent Room Type / Capacity / Facilities columns (G/H/I) to
#     match the newly assigned room, using the workbook's fixed
#     room-attribute table.
#  3. Bumps the "Generation Date" timestamp on Executive_Summary.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Basket_Course_Allocations!C9:C17 - narrowed "Allocated Rooms" lists
# ---------------------------------------------------------------------
$basketWs = $wb.Worksheets.Item("Basket_Course_Allocations")

$basketRoomChanges = @{
    9  = "C001, C002"
    10 = "C101, C102"
    11 = "C104, C205"
    12 = "C202, C203"
    13 = "C004, C204"
    14 = "C004, C102"
    15 = "C001, C104"
    16 = "C002, C202"
    17 = "C101, C203"
}

foreach ($row in $basketRoomChanges.Keys) {
    $basketWs.Range("C$row").Value = $basketRoomChanges[$row]
}

# ---------------------------------------------------------------------
# 2. Classroom_Allocation - reassigned rooms + dependent attributes
# ---------------------------------------------------------------------
$roomWs = $wb.Worksheets.Item("Classroom_Allocation")

# Fixed room -> (Room Type, Capacity, Facilities) lookup, as used
# throughout the Classroom_Allocation sheet.
$roomAttributes = @{
    "C001" = @{ RoomType = "large classroom"; Capacity = "120"; Facilities = ""                    }
    "C002" = @{ RoomType = "large classroom"; Capacity = "120"; Facilities = "Projector"            }
    "C004" = @{ RoomType = "Auditorium";      Capacity = "240"; Facilities = "Audio/Video System"   }
    "C101" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "Projector"            }
    "C102" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "Projector"            }
    "C104" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "Projector"            }
    "C202" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "Projector"            }
    "C203" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "TV"                   }
    "C204" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "TV"                   }
    "C205" = @{ RoomType = "classroom";       Capacity = "96";  Facilities = "TV"                   }
}

# New room ("M" column) per affected schedule row.
$roomReassignments = @{
    5  = "C002"; 6  = "C102"; 7  = "C205"; 8  = "C203"
    17 = "C004"; 18 = "C001"; 19 = "C002"; 20 = "C101"
    21 = "C002"; 22 = "C102"; 23 = "C205"; 24 = "C203"
    37 = "C002"; 38 = "C102"; 39 = "C205"
    41 = "C004"; 42 = "C004"; 43 = "C001"; 44 = "C002"
    45 = "C101"; 46 = "C004"; 47 = "C001"; 48 = "C002"; 49 = "C101"
    53 = "C001"; 54 = "C101"; 55 = "C104"; 56 = "C202"; 57 = "C204"
    65 = "C102"; 66 = "C104"
    69 = "C001"; 70 = "C101"; 71 = "C104"; 72 = "C202"; 73 = "C204"
    85 = "C001"; 86 = "C101"; 87 = "C104"; 88 = "C202"; 89 = "C204"
    90 = "C102"; 91 = "C104"
    94 = "C102"; 95 = "C104"; 96 = "C202"; 97 = "C203"
}

foreach ($row in $roomReassignments.Keys) {
    $newRoom = $roomReassignments[$row]
    $attrs = $roomAttributes[$newRoom]

    $roomWs.Range("M$row").Value = $newRoom
    $roomWs.Range("G$row").Value = $attrs.RoomType

    # Capacity is stored as text (e.g. "120") in this sheet, not a number;
    # force text formatting so the COM layer doesn't silently coerce the
    # numeric-looking string into a numeric cell, then restore the default
    # (unformatted) style so the cell's appearance is unchanged.
    $roomWs.Range("H$row").NumberFormat = "@"
    $roomWs.Range("H$row").Value = $attrs.Capacity
    $roomWs.Range("H$row").NumberFormat = "General"
    $roomWs.Range("H$row").Style = "Normal"

    if ($attrs.Facilities -eq "") {
        $roomWs.Range("I$row").Value = ""
    } else {
        $roomWs.Range("I$row").Value = $attrs.Facilities
    }
}

# ---------------------------------------------------------------------
# 3. Executive_Summary!C3 - regenerated timestamp
# ---------------------------------------------------------------------
$summaryWs = $wb.Worksheets.Item("Executive_Summary")
$summaryWs.Range("C3").Value = "2026-01-26 12:46"
